# Fixed a bug in flowDownSymbols: the data rows (A2:F25) were written in the
# wrong order. Re-write each row with its correct (row-id, count) values so
# the sheet matches the corrected row ordering while keeping headers (row 1)
# and the totals row (row 26) untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row($row, $a, $b, $c, $d, $e, $f) {
    $ws.Cells.Item($row, 1).Value = $a
    $ws.Cells.Item($row, 2).Value = $b
    $ws.Cells.Item($row, 3).Value = $c
    $ws.Cells.Item($row, 4).Value = $d
    $ws.Cells.Item($row, 5).Value = $e
    $ws.Cells.Item($row, 6).Value = $f
}

Set-Row 2  1001 18 30 75 60 72
Set-Row 3  501  9  52 30 75 45
Set-Row 4  601  9  60 67 60 42
Set-Row 5  401  9  48 67 75 45
Set-Row 6  701  3  90 45 97 15
Set-Row 7  1201 2  10 10 10 10
Set-Row 8  1202 2  10 10 10 10
Set-Row 9  201  9  30 15 45 30
Set-Row 10 801  3  67 65 52 45
Set-Row 11 1203 3  15 15 15 15
Set-Row 12 101  9  30 15 60 15
Set-Row 13 901  16 15 45 60 60
Set-Row 14 902  1  0  0  0  0
Set-Row 15 301  6  45 30 60 45
Set-Row 16 502  0  4  0  0  0
Set-Row 17 1    0  2  2  2  2
Set-Row 18 2    0  2  2  2  2
Set-Row 19 3    0  3  3  3  3
Set-Row 20 802  0  4  5  4  0
Set-Row 21 1101 0  15 30 30 0
Set-Row 22 602  0  0  4  0  9
Set-Row 23 402  0  0  4  0  0
Set-Row 24 702  0  0  0  4  0
Set-Row 25 1002 0  0  0  0  9
